# Update grades for midterm 1 - add missing Column F (Midterm 1) scores
# for students who did not yet have a score recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Formula  = "=51/60"
$ws.Range("F5").Formula  = "=61/60"
$ws.Range("F6").Formula  = "=51/60"
$ws.Range("F9").Formula  = "=40/60"
$ws.Range("F14").Formula = "=41/60"
$ws.Range("F28").Formula = "=58/60"
$ws.Range("F32").Formula = "=42/60"
$ws.Range("F35").Formula = "=40/60"

# Move the active cell/selection to F15, matching where the editor left off.
$ws.Range("F15").Select() | Out-Null
